$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an already-empty row (row 23) onto rows 17-22,
# then clear their contents so they become blank cells using the plain
# (non-dated) style, matching rows 23+ further down the sheet.
$ws.Range("A23:F23").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Range("A18:F18").PasteSpecial(-4122)
$ws.Range("A19:F19").PasteSpecial(-4122)
$ws.Range("A20:F20").PasteSpecial(-4122)
$ws.Range("A21:F21").PasteSpecial(-4122)
$ws.Range("A22:F22").PasteSpecial(-4122)
$ws.Range("A17:F22").ClearContents()

# Remove the now-empty trailing row 101.
$ws.Rows(101).Delete()

# Restore the view to what it looked like before the data was scrolled to.
[void]$ws.Range("D16").Select()

Write-Output "done"
